$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 180.39394
$ws.Range("I33").Value = 107.08
$ws.Range("K33").Value = 107.08
$ws.Range("M33").Value = 121.92

$ws.Range("H62").Value = 2778.3572
$ws.Range("J62").Value = 4499.8335
$ws.Range("L62").Value = 4499.8335
$ws.Range("N62").Value = -5747.8335

$ws.Range("H65").Value = 2778.3572
$ws.Range("J65").Value = 4499.8335
$ws.Range("L65").Value = 22499.1675
$ws.Range("N65").Value = -28739.1675

$ws.Range("H107").Value = 1622.875
$ws.Range("I107").Value = 1622.875
$ws.Range("K107").Value = 1622.875
$ws.Range("M107").Value = 297.125

$ws.Range("H141").Value = 12434.3
$ws.Range("I141").Value = 14430.375
$ws.Range("K141").Value = 43291.125
$ws.Range("M141").Value = -38111.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 9896802
$ws.Range("I63").Value = 17315890
$ws.Range("J63").Value = 4683.3335
$ws.Range("K63").Value = 17315890
$ws.Range("L63").Value = 4683.3335
$ws.Range("M63").Value = -17315204
$ws.Range("N63").Value = -6055.3335

$ws.Range("H66").Value = 9896802
$ws.Range("I66").Value = 17315890
$ws.Range("J66").Value = 4683.3335
$ws.Range("K66").Value = 86579450
$ws.Range("L66").Value = 23416.6675
$ws.Range("M66").Value = -86576018
$ws.Range("N66").Value = -30280.6675

$ws.Range("H74").Value = 3675.611
$ws.Range("I74").Value = 3624.1667
$ws.Range("J74").Value = 3932.8333
$ws.Range("K74").Value = 3624.1667
$ws.Range("L74").Value = 3932.8333
$ws.Range("M74").Value = -2750.1667
$ws.Range("N74").Value = -5680.8333

$ws.Range("H77").Value = 3675.611
$ws.Range("I77").Value = 3624.1667
$ws.Range("J77").Value = 3932.8333
$ws.Range("K77").Value = 18120.8335
$ws.Range("L77").Value = 19664.1665
$ws.Range("M77").Value = -13752.8335
$ws.Range("N77").Value = -28400.1665

$ws.Range("H92").Value = 38000
$ws.Range("J92").Value = 38000
$ws.Range("L92").Value = 38000
$ws.Range("N92").Value = -42992

$ws.Range("H110").Value = 837.8889
$ws.Range("I110").Value = 715.1667
$ws.Range("K110").Value = 715.1667
$ws.Range("M110").Value = 1329.8333

$ws.Range("H122").Value = 2478.7368
$ws.Range("I122").Value = 1411.091
$ws.Range("J122").Value = 3946.75
$ws.Range("K122").Value = 4233.272999999999
$ws.Range("L122").Value = 11840.25
$ws.Range("M122").Value = -1783.272999999999
$ws.Range("N122").Value = -16740.25

$ws.Range("H132").Value = 2395.102
$ws.Range("I132").Value = 1672.7106
$ws.Range("J132").Value = 4890.636
$ws.Range("K132").Value = 5018.1318
$ws.Range("L132").Value = 14671.908
$ws.Range("M132").Value = -2488.1318
$ws.Range("N132").Value = -19731.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 191.15384
$ws.Range("J80").Value = 240.11765
$ws.Range("L80").Value = 240.11765
$ws.Range("N80").Value = -2236.11765

$ws.Range("H83").Value = 191.15384
$ws.Range("J83").Value = 240.11765
$ws.Range("L83").Value = 1200.58825
$ws.Range("N83").Value = -11184.58825

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

$ws.Range("H134").Value = 2544.4927
$ws.Range("I134").Value = 1493.7646
$ws.Range("J134").Value = 3565.2
$ws.Range("K134").Value = 4481.293799999999
$ws.Range("L134").Value = 10695.6
$ws.Range("M134").Value = -1946.293799999999
$ws.Range("N134").Value = -15765.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11768477
$ws.Range("I99").Value = 25002366
$ws.Range("J99").Value = 5021
$ws.Range("K99").Value = 25002366
$ws.Range("L99").Value = 5021
$ws.Range("M99").Value = -25000868
$ws.Range("N99").Value = -8017

$ws.Range("H105").Value = 2700.1177
$ws.Range("I105").Value = 2873.818
$ws.Range("J105").Value = 2381.6667
$ws.Range("K105").Value = 2873.818
$ws.Range("L105").Value = 2381.6667
$ws.Range("M105").Value = -1126.818
$ws.Range("N105").Value = -5875.6667

$ws.Range("H107").Value = 929.6
$ws.Range("J107").Value = 1416.3334
$ws.Range("L107").Value = 1416.3334
$ws.Range("N107").Value = -5256.3334

$ws.Range("H126").Value = 11768477
$ws.Range("I126").Value = 25002366
$ws.Range("J126").Value = 5021
$ws.Range("K126").Value = 75007098
$ws.Range("L126").Value = 15063
$ws.Range("M126").Value = -75004628
$ws.Range("N126").Value = -20003

$ws.Range("H132").Value = 2515.95
$ws.Range("I132").Value = 1033.0435
$ws.Range("J132").Value = 4522.2354
$ws.Range("K132").Value = 3099.1305
$ws.Range("L132").Value = 13566.7062
$ws.Range("M132").Value = -569.1305000000002
$ws.Range("N132").Value = -18626.7062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 361
$ws.Range("I18").Value = 254.5
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 763.5
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -594.5
$ws.Range("N18").Value = -3338

$ws.Range("H113").Value = 702.3570999999999
$ws.Range("I113").Value = 616.02856
$ws.Range("J113").Value = 846.2381
$ws.Range("K113").Value = 1848.08568
$ws.Range("L113").Value = 2538.7143
$ws.Range("M113").Value = 321.9143200000001
$ws.Range("N113").Value = -6878.7143

$ws.Range("H114").Value = 2933.3684
$ws.Range("I114").Value = 162.25
$ws.Range("J114").Value = 4948.727
$ws.Range("K114").Value = 486.75
$ws.Range("L114").Value = 14846.181
$ws.Range("M114").Value = 2767.25
$ws.Range("N114").Value = -21354.181

$ws.Range("H131").Value = 6579851.5
$ws.Range("J131").Value = 816.17145
$ws.Range("L131").Value = 2448.51435
$ws.Range("N131").Value = -12528.51435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 26666.666
$ws.Range("J27").Value = 26666.666
$ws.Range("L27").Value = 26666.666
$ws.Range("N27").Value = -26998.666

$ws.Range("H107").Value = 9259990
$ws.Range("I107").Value = 276
$ws.Range("J107").Value = 11111933
$ws.Range("K107").Value = 276
$ws.Range("L107").Value = 11111933
$ws.Range("M107").Value = 1644
$ws.Range("N107").Value = -11115773

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 813.6070999999999
$ws.Range("I68").Value = 689.84906
$ws.Range("K68").Value = 689.84906
$ws.Range("M68").Value = 59.15093999999999

$ws.Range("H71").Value = 813.6070999999999
$ws.Range("I71").Value = 689.84906
$ws.Range("K71").Value = 3449.2453
$ws.Range("M71").Value = 294.7547

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 709.2778
$ws.Range("I107").Value = 544.26666
$ws.Range("J107").Value = 1534.3334
$ws.Range("K107").Value = 1632.79998
$ws.Range("L107").Value = 4603.0002
$ws.Range("M107").Value = 287.20002
$ws.Range("N107").Value = -8443.0002

$ws.Range("H132").Value = 20836084
$ws.Range("I132").Value = 1288
$ws.Range("J132").Value = 41670880
$ws.Range("K132").Value = 3864
$ws.Range("L132").Value = 125012640
$ws.Range("M132").Value = -1334
$ws.Range("N132").Value = -125017700

$ws.Range("H136").Value = 1414.5
$ws.Range("I136").Value = 591.5625
$ws.Range("J136").Value = 2355
$ws.Range("K136").Value = 1774.6875
$ws.Range("L136").Value = 7065
$ws.Range("M136").Value = 775.3125
$ws.Range("N136").Value = -12165
